$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 226.95238
$ws.Range("I28").Value = 262.35294
$ws.Range("J28").Value = 76.5
$ws.Range("K28").Value = 262.35294
$ws.Range("L28").Value = 76.5
$ws.Range("M28").Value = 222.64706
$ws.Range("N28").Value = -1046.5
$ws.Range("H33").Value = 709
$ws.Range("I33").Value = 394
$ws.Range("J33").Value = 1234
$ws.Range("K33").Value = 394
$ws.Range("L33").Value = 1234
$ws.Range("M33").Value = -165
$ws.Range("N33").Value = -1692
$ws.Range("H40").Value = 1678.9688
$ws.Range("I40").Value = 1559
$ws.Range("J40").Value = 1908
$ws.Range("K40").Value = 1559
$ws.Range("L40").Value = 1908
$ws.Range("M40").Value = -1384
$ws.Range("N40").Value = -2258
$ws.Range("H43").Value = 1786.0667
$ws.Range("I43").Value = 1250
$ws.Range("J43").Value = 1981
$ws.Range("K43").Value = 1250
$ws.Range("L43").Value = 1981
$ws.Range("M43").Value = -1181
$ws.Range("N43").Value = -2119
$ws.Range("H116").Value = 4261.696
$ws.Range("I116").Value = 3854.5454
$ws.Range("J116").Value = 4634.9165
$ws.Range("K116").Value = 3854.5454
$ws.Range("L116").Value = 4634.9165
$ws.Range("M116").Value = -412.5454
$ws.Range("N116").Value = -11518.9165
$ws.Range("H132").Value = 6253272
$ws.Range("I132").Value = 7694335.5
$ws.Range("J132").Value = 8663.333000000001
$ws.Range("K132").Value = 23083006.5
$ws.Range("L132").Value = 25989.999
$ws.Range("M132").Value = -23080476.5
$ws.Range("N132").Value = -31049.999
$ws.Range("H141").Value = 719123.4399999999
$ws.Range("I141").Value = 1599.25
$ws.Range("K141").Value = 4797.75
$ws.Range("M141").Value = 382.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4033292
$ws.Range("I2").Value = 6757474.5
$ws.Range("J2").Value = 1501.8
$ws.Range("K2").Value = 6757474.5
$ws.Range("L2").Value = 1501.8
$ws.Range("M2").Value = -6757361.5
$ws.Range("N2").Value = -1727.8
$ws.Range("H5").Value = 62.42857
$ws.Range("I5").Value = 56.166668
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 56.166668
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 55.833332
$ws.Range("N5").Value = -324
$ws.Range("H116").Value = 4033292
$ws.Range("I116").Value = 6757474.5
$ws.Range("J116").Value = 1501.8
$ws.Range("K116").Value = 6757474.5
$ws.Range("L116").Value = 1501.8
$ws.Range("M116").Value = -6755180.5
$ws.Range("N116").Value = -6089.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4033292
$ws.Range("I3").Value = 6757474.5
$ws.Range("J3").Value = 1501.8
$ws.Range("K3").Value = 6757474.5
$ws.Range("L3").Value = 1501.8
$ws.Range("M3").Value = -6757360.5
$ws.Range("N3").Value = -1729.8
$ws.Range("H4").Value = 62.42857
$ws.Range("I4").Value = 56.166668
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 56.166668
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 58.833332
$ws.Range("N4").Value = -330

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M15").ClearContents()
$ws.Range("M51").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("H15").Value = 2400
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 2400
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2400
$ws.Range("N15").Value = -2740
$ws.Range("H22").Value = 740.08
$ws.Range("I22").Value = 303.125
$ws.Range("J22").Value = 1516.8889
$ws.Range("K22").Value = 303.125
$ws.Range("L22").Value = 1516.8889
$ws.Range("M22").Value = 46.875
$ws.Range("N22").Value = -2216.8889
$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 20000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21472
$ws.Range("H56").Value = 10001
$ws.Range("I56").Value = 10001
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 10001
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -9156
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 20000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 20000
$ws.Range("N61").Value = -20696
$ws.Range("H62").Value = 3345.4546
$ws.Range("I62").Value = 2420.2222
$ws.Range("K62").Value = 2420.2222
$ws.Range("M62").Value = -1796.2222
$ws.Range("H65").Value = 3345.4546
$ws.Range("I65").Value = 2420.2222
$ws.Range("K65").Value = 12101.111
$ws.Range("M65").Value = -8981.111000000001
$ws.Range("H107").Value = 1149.0358
$ws.Range("I107").Value = 1258.0769
$ws.Range("J107").Value = 1054.5333
$ws.Range("K107").Value = 1258.0769
$ws.Range("L107").Value = 1054.5333
$ws.Range("M107").Value = 661.9231
$ws.Range("N107").Value = -4894.5333
$ws.Range("H132").Value = 2890.8572
$ws.Range("I132").Value = 2372.7222
$ws.Range("J132").Value = 5999.6665
$ws.Range("K132").Value = 7118.1666
$ws.Range("L132").Value = 17998.9995
$ws.Range("M132").Value = -4588.1666
$ws.Range("N132").Value = -23058.9995
$ws.Range("H141").Value = 21549.834
$ws.Range("J141").Value = 23314.285
$ws.Range("L141").Value = 23314.285
$ws.Range("N141").Value = -33674.285

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2579.8
$ws.Range("J57").Value = 4000
$ws.Range("L57").Value = 12000
$ws.Range("N57").Value = -13118
$ws.Range("H59").Value = 2843.7144
$ws.Range("I59").Value = 1200
$ws.Range("K59").Value = 3600
$ws.Range("M59").Value = -3060
$ws.Range("H63").Value = 13433.111
$ws.Range("I63").Value = 1779.6
$ws.Range("J63").Value = 28000
$ws.Range("K63").Value = 5338.799999999999
$ws.Range("L63").Value = 84000
$ws.Range("M63").Value = -4589.799999999999
$ws.Range("N63").Value = -85498
$ws.Range("H64").Value = 8175.0586
$ws.Range("I64").Value = 1496.3334
$ws.Range("K64").Value = 4489.0002
$ws.Range("M64").Value = -4219.0002
$ws.Range("H66").Value = 13433.111
$ws.Range("I66").Value = 1779.6
$ws.Range("J66").Value = 28000
$ws.Range("K66").Value = 16016.4
$ws.Range("L66").Value = 252000
$ws.Range("M66").Value = -12272.4
$ws.Range("N66").Value = -259488
$ws.Range("H67").Value = 8175.0586
$ws.Range("I67").Value = 1496.3334
$ws.Range("K67").Value = 4489.0002
$ws.Range("M67").Value = -3553.0002
$ws.Range("H118").Value = 2348.8948
$ws.Range("I118").Value = 376.33334
$ws.Range("J118").Value = 2718.75
$ws.Range("K118").Value = 1129.00002
$ws.Range("L118").Value = 8156.25
$ws.Range("M118").Value = 113.9999800000001
$ws.Range("N118").Value = -10642.25
$ws.Range("H131").Value = 1156.71
$ws.Range("I131").Value = 3464.1428
$ws.Range("J131").Value = 983.0323
$ws.Range("K131").Value = 10392.4284
$ws.Range("L131").Value = 2949.0969
$ws.Range("M131").Value = -5352.428400000001
$ws.Range("N131").Value = -13029.0969

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3480.1538
$ws.Range("I80").Value = 3019.7334
$ws.Range("J80").Value = 4108
$ws.Range("K80").Value = 3019.7334
$ws.Range("L80").Value = 4108
$ws.Range("M80").Value = -2021.7334
$ws.Range("N80").Value = -6104
$ws.Range("H83").Value = 3480.1538
$ws.Range("I83").Value = 3019.7334
$ws.Range("J83").Value = 4108
$ws.Range("K83").Value = 15098.667
$ws.Range("L83").Value = 20540
$ws.Range("M83").Value = -10106.667
$ws.Range("N83").Value = -30524
$ws.Range("H102").Value = 54264.4
$ws.Range("I102").Value = 3338.5
$ws.Range("J102").Value = 105190.3
$ws.Range("K102").Value = 3338.5
$ws.Range("L102").Value = 105190.3
$ws.Range("M102").Value = -1716.5
$ws.Range("N102").Value = -108434.3
$ws.Range("H113").Value = 3212.6553
$ws.Range("I113").Value = 2616.682
$ws.Range("J113").Value = 5085.7144
$ws.Range("K113").Value = 2616.682
$ws.Range("L113").Value = 5085.7144
$ws.Range("M113").Value = -446.6819999999998
$ws.Range("N113").Value = -9425.714400000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 8166.6665
$ws.Range("J3").Value = 8166.6665
$ws.Range("L3").Value = 8166.6665
$ws.Range("N3").Value = -8390.666499999999
$ws.Range("H14").Value = 219587.42
$ws.Range("J14").Value = 12504.5
$ws.Range("L14").Value = 12504.5
$ws.Range("N14").Value = -12848.5
$ws.Range("H15").Value = 8166.6665
$ws.Range("J15").Value = 8166.6665
$ws.Range("L15").Value = 8166.6665
$ws.Range("N15").Value = -8506.666499999999
$ws.Range("H46").Value = 1271.2
$ws.Range("I46").Value = 536
$ws.Range("J46").Value = 1455
$ws.Range("K46").Value = 536
$ws.Range("L46").Value = 1455
$ws.Range("M46").Value = -348
$ws.Range("N46").Value = -1831
$ws.Range("H68").Value = 2763.5881
$ws.Range("I68").Value = 1170.7273
$ws.Range("J68").Value = 5683.8335
$ws.Range("K68").Value = 1170.7273
$ws.Range("L68").Value = 5683.8335
$ws.Range("M68").Value = -421.7273
$ws.Range("N68").Value = -7181.8335
$ws.Range("H71").Value = 2763.5881
$ws.Range("I71").Value = 1170.7273
$ws.Range("J71").Value = 5683.8335
$ws.Range("K71").Value = 5853.636500000001
$ws.Range("L71").Value = 28419.1675
$ws.Range("M71").Value = -2109.636500000001
$ws.Range("N71").Value = -35907.1675
$ws.Range("H100").Value = 2810
$ws.Range("I100").Value = 1766.6666
$ws.Range("J100").Value = 3436
$ws.Range("K100").Value = 1766.6666
$ws.Range("L100").Value = 3436
$ws.Range("M100").Value = -1225.6666
$ws.Range("N100").Value = -4518
